$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(255, 1).Value = 'z0015'
$ws.Cells.Item(255, 2).Value = '포물선의 초점의 좌표를 구합니다.'
$ws.Cells.Item(255, 3).Value = '32111_z28'

$ws.Cells.Item(256, 1).Value = 'z0016'
$ws.Cells.Item(256, 2).Value = '주어진 두 점 사이의 거리를 이용해서 방정식을 세웁니다.'

$ws.Cells.Item(257, 1).Value = 'z0017'
$ws.Cells.Item(257, 2).Value = '\mathrm{PQ}$를 빗변으로 하는 직각삼각형을 이용해서 방정식을 세웁니다.'

$ws.Cells.Item(258, 1).Value = 'z0018'
$ws.Cells.Item(258, 2).Value = '\mathrm{P}$, \mathrm{Q}$의 $x$좌표를 미지수로 도입해서 방정식을 풉니다.'

$ws.Cells.Item(259, 1).Value = 'z0019'
$ws.Cells.Item(259, 2).Value = '조건 (가)를 만족시키는 점 $\mathrm{P}$의 영역을 구합니다.'
$ws.Cells.Item(259, 3).Value = '32111_z29'

$ws.Cells.Item(260, 1).Value = 'z0020'
$ws.Cells.Item(260, 2).Value = '조건 (나)를 벡터의 시점 $\mathrm{O}$에 대해 정리하고 주어진 길이와 각의 조건을 이용해 식을 정리합니다.'

$ws.Cells.Item(261, 1).Value = 'z0021'
$ws.Cells.Item(261, 2).Value = '$|3 \overrightarrow{\mathrm{OP}}-\overrightarrow{\mathrm{OX}}|$가 최대가 되도록 하는 경우를 발견해 그 최댓값을 구합니다.'

$ws.Cells.Item(262, 1).Value = 'z0022'
$ws.Cells.Item(262, 2).Value = '$|3 \overrightarrow{\mathrm{OP}}-\overrightarrow{\mathrm{OX}}|$가 최소가 되도록 하는 경우를 발견해 그 최솟값을 구합니다.'

$ws.Cells.Item(263, 1).Value = 'z0023'
$ws.Cells.Item(263, 2).Value = '정사영된 삼각형의 각 꼭짓점을 파악해서 정사영된 삼각형 넓이의 최댓값을 구합니다.'
$ws.Cells.Item(263, 3).Value = '32111_z30'

$ws.Cells.Item(264, 1).Value = 'z0024'
$ws.Cells.Item(264, 2).Value = '피타고라스 정리를 이용해서 정사영의 넓이를 최대화 시키는 삼각형 $\mathrm{PQR}$의 세변의 길이를 구합니다.'

$ws.Cells.Item(265, 1).Value = 'z0025'
$ws.Cells.Item(265, 2).Value = '삼각형 $\mathrm{PQR}$의 넓이를 구합니다.'

$ws.Cells.Item(266, 1).Value = 'z0026'
$ws.Cells.Item(266, 2).Value = '정사영 전과 후의 두 넓이의 비율을 이용해서 두 평면이 이루는 예각에 대한 코사인 값을 구합니다. '

$ws.Cells.Item(267, 1).Value = 'z0027'
$ws.Cells.Item(267, 2).Value = '코사인 값을 이용해서 정사영의 넓이를 구합니다.'

$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 241
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B269").Select()
